# "Generate Report for Handback" -- the df2e5e7c*.md file has now been
# handed back in sync, so it drops off the report; the still-open
# 6ae3d7f8*.md file picks up newer handoff/handback timestamps.

$wb = $excel.ActiveWorkbook
$hyperlinkColor = 15570276   # BGR long for RGB 6495ED (CornflowerBlue), matches the workbook's existing HyperLink style

function Restyle-Hyperlinks($ws, [string[]]$addrs) {
    foreach ($addr in $addrs) {
        $c = $ws.Range($addr)
        $c.Font.Underline = $true
        $c.Font.Color = $hyperlinkColor
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn": row 2 (6ae3d7f8...) gets refreshed handoff/handback
# timestamps; row 3 (df2e5e7c...) is removed entirely.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-12 18:39:40"
$wsZh.Range("H2").Value = "2016-03-12 18:39:56"

$wsZh.Hyperlinks.Delete()
$wsZh.Rows.Item(3).Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/84e45eb03401a08ba8aadd6d7813bd115f6133de/e2e/6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md", [Type]::Missing, [Type]::Missing, "6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/84e45eb03401a08ba8aadd6d7813bd115f6133de/e2e/6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d0d53f0ff09089c866d2dfa88a011dc63a71647c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6ae3d7f8-b5c4-4c63-a241-d93bd81600db.b69b922b01c1db1121585ad9a66e90bb3efc8080.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "6ae3d7f8-b5c4-4c63-a241-d93bd81600db.b69b922b01c1db1121585ad9a66e90bb3efc8080.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/9cb8f9eb32c18d971caf38001043dea239129ae6/e2e/6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md", [Type]::Missing, [Type]::Missing, "6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9db018241bf0bda6b4daf54370c1248aa0b831f7/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6ae3d7f8-b5c4-4c63-a241-d93bd81600db.b69b922b01c1db1121585ad9a66e90bb3efc8080.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "6ae3d7f8-b5c4-4c63-a241-d93bd81600db.b69b922b01c1db1121585ad9a66e90bb3efc8080.zh-cn.xlf") | Out-Null

Restyle-Hyperlinks $wsZh @("A2","B2","D2","F2","G2")

# ---------------------------------------------------------------------
# Sheet "de-de": same treatment, German handoff/handback timestamps.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-12 18:39:43"
$wsDe.Range("H2").Value = "2016-03-12 18:40:03"

$wsDe.Hyperlinks.Delete()
$wsDe.Rows.Item(3).Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/84e45eb03401a08ba8aadd6d7813bd115f6133de/e2e/6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md", [Type]::Missing, [Type]::Missing, "6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/84e45eb03401a08ba8aadd6d7813bd115f6133de/e2e/6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3949fca6c63c6023c3c9332e8ce85da4d750897e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6ae3d7f8-b5c4-4c63-a241-d93bd81600db.b69b922b01c1db1121585ad9a66e90bb3efc8080.de-de.xlf", [Type]::Missing, [Type]::Missing, "6ae3d7f8-b5c4-4c63-a241-d93bd81600db.b69b922b01c1db1121585ad9a66e90bb3efc8080.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/63b7124ee0fb0675314ceabd58be9cabd0e31383/e2e/6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md", [Type]::Missing, [Type]::Missing, "6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3e909c14a5c0d2546d1eb7e61d78507efa55fd69/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6ae3d7f8-b5c4-4c63-a241-d93bd81600db.b69b922b01c1db1121585ad9a66e90bb3efc8080.de-de.xlf", [Type]::Missing, [Type]::Missing, "6ae3d7f8-b5c4-4c63-a241-d93bd81600db.b69b922b01c1db1121585ad9a66e90bb3efc8080.de-de.xlf") | Out-Null

Restyle-Hyperlinks $wsDe @("A2","B2","D2","F2","G2")

# ---------------------------------------------------------------------
# Sheet "Overview": drop the df2e5e7c... summary row (row 3), keep the
# still-pending 6ae3d7f8... row.
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Hyperlinks.Delete()
$wsOv.Rows.Item(3).Delete()

$wsOv.Hyperlinks.Add($wsOv.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/84e45eb03401a08ba8aadd6d7813bd115f6133de/e2e/6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md", [Type]::Missing, [Type]::Missing, "6ae3d7f8-b5c4-4c63-a241-d93bd81600db.md") | Out-Null

Restyle-Hyperlinks $wsOv @("A2")
